$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the sample data values in row 2 (test case data refresh)
$ws.Range("A2").Value = "Keeru"
$ws.Range("B2").Value = "keerthi"
$ws.Range("C2").Value = "keerthi232@gmail.com"
$ws.Range("D2").Value = "keerthi@232"
$ws.Range("E2").Value = "keerthi@232"

# Widen column B to fit the new content
$ws.Columns.Item(2).ColumnWidth = 13.29
